# The "Drop Downs" sheet used to hold 4 side-by-side lists (with a bold
# header row in row 1): Platform/Device Types, Ad Format Types,
# Yes/No Options, Pricing Models.
#
# The new layout drops the header row entirely and stacks each list in a
# single column, one after another, starting at row 3:
#   A3:A8   -> Platform/Device Types values
#   B9:B16  -> Ad Format Types values
#   C17:C18 -> Yes/No Options values
#   D19:D28 -> Pricing Models values

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Drop Downs")

# Wipe the existing grid (values + the bold header formatting) so the sheet
# can be rebuilt from scratch with the new stacked layout.
$ws.Cells.Clear()

$platformDeviceTypes = @(
    "Desktop",
    "Mobile - Web",
    "Mobile - In App",
    "Tablet",
    "Connected TV",
    "Cross-Platform"
)

$adFormatTypes = @(
    "Rich Media - W/O Video",
    "Rich Media - WITH Video",
    "Tracking - 1x1",
    "Site Served",
    "Standard Banner",
    "HTML5 Standard Banner",
    "VAST Video",
    "In-Stream Video"
)

$yesNoOptions = @(
    "Yes",
    "No"
)

$pricingModels = @(
    "CPM",
    "Flat Fee",
    "Added Value",
    "CPC",
    "CPE",
    "CPV",
    "CPCV",
    "CPA",
    "vCPM",
    "dCPM"
)

$row = 3
foreach ($val in $platformDeviceTypes) {
    $ws.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}

$row = 9
foreach ($val in $adFormatTypes) {
    $ws.Cells.Item($row, 2).Value = $val
    $row = $row + 1
}

$row = 17
foreach ($val in $yesNoOptions) {
    $ws.Cells.Item($row, 3).Value = $val
    $row = $row + 1
}

$row = 19
foreach ($val in $pricingModels) {
    $ws.Cells.Item($row, 4).Value = $val
    $row = $row + 1
}
